$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Modyifiying with some changes"
$ws.Range("A2").Value = "New rows are added to the column"

$ws.Range("A3").Select() | Out-Null
